# Insert a new data row at row 82 (pushing existing rows 82..98 down to 83..99)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(82).Insert()

$ws.Cells.Item(82, 1).Value = 6
$ws.Cells.Item(82, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(82, 3).Value = "Metropolitana"
$ws.Cells.Item(82, 4).Value = 45127
$ws.Cells.Item(82, 5).Value = 13
$ws.Cells.Item(82, 6).Value = 100112035
$ws.Cells.Item(82, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 200
$ws.Cells.Item(82, 11).Value = 17000
$ws.Cells.Item(82, 12).Value = 18000
$ws.Cells.Item(82, 13).Value = 17400
$ws.Cells.Item(82, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(82, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(82, 16).Value = 1160
$ws.Cells.Item(82, 17).Value = 15
$ws.Cells.Item(82, 18).Value = "Hortaliza"
